$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (16th May Refresh) - appended after the existing last row (33)
$newRows = @(
    @(10005, 110033, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110034, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110035, 10005, "eng", $true, "superadmin", "now()")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Match the selection/view state left behind after entering the new rows
$ws.Range("A37:XFD1048576").Select()
